$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q1" worksheet right before the "总计" (total) sheet.
# ---------------------------------------------------------------------------
$totalSheetBefore = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheetBefore)
$newSheet.Name = "2022-Q1"

# Match the page layout used by the other quarter sheets.
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Borrow the header / index-column formatting from an existing quarter sheet.
$srcSheet = $wb.Worksheets.Item("2021-Q3")
$srcSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$srcSheet.Range("A2").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row values.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2 - columns B..G are stored as text, A and H are numeric.
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "014125"
$newSheet.Range("C2").Value = "华夏中证1000指数增强A"
$newSheet.Range("D2").Value = "7.03"
$newSheet.Range("E2").Value = "89.75"
$newSheet.Range("F2").Value = "0.81"
$newSheet.Range("G2").Value = "0.0569"
$newSheet.Range("B2:G2").Style = "Normal"
$newSheet.Range("H2").Value = 9

# Row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3:G3").NumberFormat = "@"
$newSheet.Range("B3").Value = "014126"
$newSheet.Range("C3").Value = "华夏中证1000指数增强C"
$newSheet.Range("D3").Value = "6.09"
$newSheet.Range("E3").Value = "89.75"
$newSheet.Range("F3").Value = "0.81"
$newSheet.Range("G3").Value = "0.0493"
$newSheet.Range("B3:G3").Style = "Normal"
$newSheet.Range("H3").Value = 9

# ---------------------------------------------------------------------------
# 2. Add a "2022-Q1" summary row at the top of the "总计" sheet's data and
#    renumber the existing index column. Re-fetch the sheet by name since the
#    earlier reference's position shifted once the new sheet was inserted.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# The inserted row loses the index-column formatting - copy it back from the
# row immediately below (which still carries the original formatting).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2:D2").Style = "Normal"
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.11

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3

# Restore the originally active sheet/tab (adding a sheet shifts focus to it).
$wb.Worksheets.Item(1).Activate()
